# Replace the three embedded placeholder images with plain hyperlink runs
# that point at the real image URLs on ura.gov.sg, per the commit diff.
#
# Each image is a lone InlineShape inside its own paragraph (FirstParagraph
# or BodyText). We delete the shape and insert, in its place, a hyperlink
# run styled with the built-in "Hyperlink" character style whose visible
# text is the image's URL.

$d = $word.ActiveDocument

function Replace-ImageWithHyperlink($descr, $url) {
    $found = $null
    for ($i = 1; $i -le $d.InlineShapes.Count; $i++) {
        $shp = $d.InlineShapes.Item($i)
        if ($shp.AlternativeText -eq $descr) {
            $found = $shp
            break
        }
    }
    if ($found -eq $null) {
        $found = $d.InlineShapes.Item(1)
    }
    $rng = $found.Range
    $rng.Collapse(1)
    $found.Delete()
    $h = $d.Hyperlinks.Add($rng, $url, "", "", $url)
}

Replace-ImageWithHyperlink "Road buffer and setback" "https://ura.gov.sg/-/media/Corporate/Guidelines/Development-control/Commercial/C05_Road_Buffer_and_Setbacks.jpg?h=100%25&w=100%25"
Replace-ImageWithHyperlink "" "https://ura.gov.sg/-/media/Corporate/Guidelines/Development-control/Commercial/C11_Setbacks_for_Ancillary_Structures_Substation.jpg?h=100%25&w=100%25"
Replace-ImageWithHyperlink "Setback for multi-storey car parks" "https://ura.gov.sg/-/media/Corporate/Guidelines/Development-control/Commercial/C06_Setback_for_MSCP.jpg?h=100%25&w=100%25"

Write-Output ("Remaining InlineShapes: " + $d.InlineShapes.Count)
Write-Output ("Hyperlinks: " + $d.Hyperlinks.Count)
